$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.924.54"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.898.19"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'0.7546"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'240.40"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.3049"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("D9").Value = "'25.46"
$ws.Range("E9").Value = "  -6.02%  "
$ws.Range("D10").Value = "'0.06848"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").Value = "'0.07978"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'0.7458"
$ws.Range("E12").Value = "  -4.62%  "
$ws.Range("D13").Value = "1.897.67"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'5.191"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "'91.35"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "29.933.44"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'13.95"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").Value = "'5.978"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").Value = "'243.74"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'0.000007725"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'6.973"
$ws.Range("E23").Value = "  +4.05%  "
$ws.Range("D24").Value = "'9.251"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").Value = "'165.41"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'18.73"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").Value = "'0.1278"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'2.027"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "'1.389"
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("D30").Value = "'1.517"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").Value = "'4.276"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "'4.024"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").Value = "'0.05322"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").Value = "'1.250"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("D35").Value = "'0.7239"
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").Value = "'2.717"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "'0.01913"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'2.789"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "'6.175"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "'0.4406"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "'72.12"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'1.892"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "'0.8258"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").Value = "'100.89"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "'7.533"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").Value = "'9.791"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "2.050.47"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'36.29"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "'0.05967"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "'1.472"
$ws.Range("E51").Value = "  +0.21%  "
